# Edit the "register" sheet test data and selection, as captured by the
# commit "changes made to page and testData".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("register")

# --- testData changes -----------------------------------------------
# First registrant: first name "Javad" -> "Mohammad " and a longer/updated
# gmail address.
$ws.Range("A2").Value = "Mohammad "
$ws.Range("C2").Value = "m.jawadadeli123450012222@gmail.com"

# Second registrant: first name "Azima" -> "Azimeh" and a longer/updated
# ucf.edu address.
$ws.Range("A3").Value = "Azimeh"
$ws.Range("C3").Value = "azimeh.kazemian14568633309909@ucf.edu"

# --- page (view) changes ---------------------------------------------
# Move the active selection on the register sheet to E9.
$ws.Range("E9").Select() | Out-Null

# Resize the workbook window (best effort; some hosts fix the window
# geometry from the original file and ignore this).
$excel.ActiveWindow.Width = 13665
